# Apply "requirement separation" edit to the RELS worksheet.
#
# Original layout: A=Course name, B=Units, C=Prerequisites,
#                   D=Terms Typically Offered
# New layout:       A=Course name, B=Units, C=Prerequisites,
#                   D=Corequisites, E=Concurrent, F=Recommended,
#                   G=Terms Typically Offered
#
# The old column D (term data) is pushed three columns to the right by
# inserting three new, blank columns at D. We then populate the new
# D/E/F columns with header + "NA" placeholder data, and fix up row 18
# whose prerequisite text embedded a "Recommended:" clause that now
# belongs in its own column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing "Terms Typically Offered" column (D) three slots to
# the right, to G, leaving three empty columns (D:F) behind.
$ws.Columns("D:F").Insert()

# New header row.
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"
$ws.Range("G1").Value = "Terms Typically Offered"

# Default every data row (2-20) in the three new columns to "NA".
for ($r = 2; $r -le 20; $r++) {
    $ws.Cells.Item($r, 4).Value = "NA"
    $ws.Cells.Item($r, 5).Value = "NA"
    $ws.Cells.Item($r, 6).Value = "NA"
}

# Row 18 (RELS 380) previously had its "Recommended" clause folded into
# the Prerequisites text; split it out into its own column and trim the
# Prerequisites cell.
$ws.Range("C18").Value = "Completion of GE Area A with grades of C- or better."
$ws.Range("F18").Value = "Completion of one class in POLS or RELS."
$ws.Range("G18").Value = "W "
